# Add a new "writeReview" worksheet right before "signIn" (so the sheet
# order becomes: offlineShopingProcess, addToCartAndVerifyPrice,
# writeReview, signIn) and populate it with the review test data.

$wb = $excel.ActiveWorkbook

$signIn = $wb.Worksheets.Item("signIn")
$newSheet = $wb.Worksheets.Add($signIn)
$newSheet.Name = "writeReview"

# Header row
$newSheet.Range("A1").Value = "itemname"
$newSheet.Range("B1").Value = "nick"
$newSheet.Range("C1").Value = "summary"
$newSheet.Range("D1").Value = "content"

# Data row
$newSheet.Range("A2").Value = "Josie Yoga Jacket"
$newSheet.Range("B2").Value = "janusz"
$newSheet.Range("C2").Value = "niezła"
$newSheet.Range("D2").Value = "niezła bluzka, taka niezbyt wygodna"

$newSheet.PageSetup.Orientation = 1

# Make the new sheet the active one, with D2 selected, mirroring the
# authored workbook (tabSelected moves from addToCartAndVerifyPrice to
# writeReview, activeTab becomes the new sheet's index).
[void]$newSheet.Activate()
[void]$newSheet.Range("D2").Select()
